$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => column letter => new text value, extracted from the authoritative diff.
$updates = @{
    2 = @{ D = "65.926.89"; E = "  +1.14%  " }
    3 = @{ D = "3.202.35"; E = "  +0.82%  " }
    4 = @{ E = "  +0.04%  " }
    5 = @{ D = "601.16"; E = "  +3.62%  " }
    6 = @{ D = "153.32"; E = "  +0.99%  " }
    7 = @{ E = "  +0.08%  " }
    8 = @{ D = "3.200.72"; E = "  +0.81%  " }
    9 = @{ D = "0.532"; E = "  +0.11%  " }
    10 = @{ D = "0.159"; E = "  -1.98%  " }
    11 = @{ D = "6.10"; E = "  -1.79%  " }
    12 = @{ D = "0.510"; E = "  +1.15%  " }
    13 = @{ D = "0.0000270"; E = "  -0.98%  " }
    14 = @{ D = "39.36"; E = "  +4.39%  " }
    15 = @{ D = "3.728.98"; E = "  +0.93%  " }
    16 = @{ D = "7.48"; E = "  +4.00%  " }
    17 = @{ D = "66.006.22"; E = "  +1.20%  " }
    18 = @{ D = "3.209.51"; E = "  +1.13%  " }
    19 = @{ E = "  +0.00%  " }
    20 = @{ D = "510.32"; E = "  -0.68%  " }
    21 = @{ D = "15.40"; E = "  +3.33%  " }
    22 = @{ D = "0.739"; E = "  +1.36%  " }
    23 = @{ D = "8.14"; E = "  +3.56%  " }
    24 = @{ D = "15.37"; E = "  -0.77%  " }
    25 = @{ D = "84.89"; E = "  -0.41%  " }
    26 = @{ E = "  -0.21%  " }
    27 = @{ E = "  +1.91%  " }
    28 = @{ D = "3.01"; E = "  +2.40%  " }
    29 = @{ D = "2.27"; E = "  +3.05%  " }
    30 = @{ D = "2.88"; E = "  +1.42%  " }
    31 = @{ D = "6.87"; E = "  +8.45%  " }
    32 = @{ D = "28.04"; E = "  +0.39%  " }
    33 = @{ E = "  +1.90%  " }
    34 = @{ E = "  +0.25%  " }
    35 = @{ D = "6.57"; E = "  -0.73%  " }
    36 = @{ D = "55.00"; E = "  -1.61%  " }
    37 = @{ D = "0.0904"; E = "  -0.14%  " }
    38 = @{ D = "486.18"; E = "  +1.62%  " }
    39 = @{ E = "  -0.80%  " }
    40 = @{ D = "2.94"; E = "  -5.64%  " }
    41 = @{ D = "8.90"; E = "  +2.31%  " }
    42 = @{ D = "0.303"; E = "  +6.07%  " }
    43 = @{ D = "0.120"; E = "  +1.66%  " }
    44 = @{ D = "2.953.88"; E = "  -4.27%  " }
    45 = @{ E = "  +6.14%  " }
    46 = @{ E = "  -0.93%  " }
    47 = @{ D = "28.54"; E = "  -3.05%  " }
    48 = @{ E = "  +0.06%  " }
    49 = @{ E = "  +0.58%  " }
    50 = @{ D = "2.31"; E = "  +1.90%  " }
    51 = @{ D = "120.32"; E = "  -0.23%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    if ($vals.ContainsKey("D")) {
        # Column D cells hold text that LOOKS numeric ("6.10", "0.0000270", ...).
        # Force the cell to Text format first so Excel keeps the exact digits/
        # trailing zeros instead of silently re-parsing the string as a Number,
        # then clear the formatting delta so no stray style is left behind.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $vals["D"]
        $dCell.ClearFormats()
    }

    if ($vals.ContainsKey("E")) {
        # Column E values ("  +1.14%  ") already contain spaces/% so Excel
        # keeps them as text without any extra coercion needed.
        $ws.Cells.Item($row, 5).Value = $vals["E"]
    }
}
